# Apply updated "想去人数" (F) / "最低票价" (G) figures to the
# 广州-漫展信息 workbook, reflecting the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> column -> new value
$changes = @{
    "展览" = @{
        2  = @{ F = 13793 }
        6  = @{ F = 809 }
        7  = @{ F = 2204 }
        8  = @{ F = 219 }
        9  = @{ F = 141 }
        10 = @{ F = 125 }
        11 = @{ F = 259; G = 60 }
        13 = @{ F = 630 }
        14 = @{ F = 472 }
        15 = @{ F = 549 }
        16 = @{ F = 350 }
        17 = @{ F = 40 }
        18 = @{ F = 324 }
        19 = @{ F = 915 }
        20 = @{ F = 172 }
        21 = @{ F = 98 }
        22 = @{ F = 55 }
        25 = @{ F = 135 }
        26 = @{ F = 51 }
    }
    "演出" = @{
        4  = @{ F = 147 }
        6  = @{ F = 147 }
        8  = @{ F = 2201 }
        15 = @{ F = 1948 }
    }
    "本地生活" = @{
        3 = @{ F = 229 }
        4 = @{ F = 142 }
    }
    "全部类型" = @{
        3  = @{ F = 13793 }
        7  = @{ F = 809 }
        10 = @{ F = 2204 }
        11 = @{ F = 229 }
        12 = @{ F = 219 }
        13 = @{ F = 141 }
        14 = @{ F = 125 }
        15 = @{ F = 259; G = 60 }
        16 = @{ F = 147 }
        19 = @{ F = 147 }
        20 = @{ F = 142 }
        21 = @{ F = 630 }
        22 = @{ F = 472 }
        23 = @{ F = 549 }
        24 = @{ F = 350 }
        25 = @{ F = 40 }
        26 = @{ F = 324 }
        27 = @{ F = 915 }
        29 = @{ F = 2201 }
        34 = @{ F = 172 }
        35 = @{ F = 98 }
        36 = @{ F = 55 }
        41 = @{ F = 135 }
        42 = @{ F = 51 }
        43 = @{ F = 1948 }
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $changes[$sheetName]
    foreach ($row in $rowMap.Keys) {
        $colMap = $rowMap[$row]
        foreach ($col in $colMap.Keys) {
            $addr = "$col$row"
            $ws.Range($addr).Value = $colMap[$col]
        }
    }
}
